# Corrected bug in budgetController and generateStatementController
# Fix: transactions that are of type "Expense" were incorrectly labeled
# as "Income" in the generated statement. Update the Type column (B)
# for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Expense"
$ws.Range("B4").Value = "Expense"
$ws.Range("B5").Value = "Expense"
